# Apply the "cryptos list" data refresh described in the commit:
#   "Updated cryptos list on Tue Dec 12 11:56:37 UTC 2023 with GitHub Actions"
#
# For every changed row, column D holds the coin Price and column E holds the
# 1h Volume percentage change; both are plain text cells (not numbers), matching
# the source data which keeps the "."-grouped price strings and padded "%"
# strings as text. Values that would otherwise be auto-recognized by Excel as a
# number (e.g. "4.04") are entered with a leading apostrophe, exactly like a
# user forcing text entry in the UI, so they stay text cells after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "41.676.58" }
    @{ Cell = "E2"; Value = "  -1.56%  " }
    @{ Cell = "D3"; Value = "2.212.94" }
    @{ Cell = "E3"; Value = "  -1.58%  " }
    @{ Cell = "E4"; Value = "  +0.08%  " }
    @{ Cell = "D5"; Value = "'251.06" }
    @{ Cell = "E5"; Value = "  +6.30%  " }
    @{ Cell = "D6"; Value = "'0.629" }
    @{ Cell = "E6"; Value = "  +0.47%  " }
    @{ Cell = "D7"; Value = "'71.03" }
    @{ Cell = "E7"; Value = "  +1.68%  " }
    @{ Cell = "E8"; Value = "  +0.12%  " }
    @{ Cell = "E9"; Value = "  +7.72%  " }
    @{ Cell = "D10"; Value = "'40.55" }
    @{ Cell = "E10"; Value = "  +10.49%  " }
    @{ Cell = "E11"; Value = "  -3.09%  " }
    @{ Cell = "D12"; Value = "'58.31" }
    @{ Cell = "E12"; Value = "  -0.87%  " }
    @{ Cell = "D13"; Value = "'7.21" }
    @{ Cell = "E13"; Value = "  +6.53%  " }
    @{ Cell = "E14"; Value = "  -0.82%  " }
    @{ Cell = "D15"; Value = "2.545.71" }
    @{ Cell = "E15"; Value = "  -1.33%  " }
    @{ Cell = "E16"; Value = "  -1.08%  " }
    @{ Cell = "D17"; Value = "'0.872" }
    @{ Cell = "E17"; Value = "  -0.58%  " }
    @{ Cell = "D18"; Value = "2.208.05" }
    @{ Cell = "E18"; Value = "  -1.77%  " }
    @{ Cell = "D19"; Value = "41.687.79" }
    @{ Cell = "E19"; Value = "  -1.32%  " }
    @{ Cell = "D20"; Value = "0.0₃0959" }
    @{ Cell = "E20"; Value = "  -1.75%  " }
    @{ Cell = "D21"; Value = "'6.21" }
    @{ Cell = "E21"; Value = "  -1.16%  " }
    @{ Cell = "D22"; Value = "'72.63" }
    @{ Cell = "E22"; Value = "  -1.10%  " }
    @{ Cell = "D23"; Value = "'234.47" }
    @{ Cell = "E23"; Value = "  -0.92%  " }
    @{ Cell = "E24"; Value = "  +2.19%  " }
    @{ Cell = "D25"; Value = "'4.04" }
    @{ Cell = "E25"; Value = "  +10.77%  " }
    @{ Cell = "E26"; Value = "  -0.05%  " }
    @{ Cell = "E27"; Value = "  +4.72%  " }
    @{ Cell = "D28"; Value = "'11.03" }
    @{ Cell = "E28"; Value = "  +9.96%  " }
    @{ Cell = "E29"; Value = "  -2.41%  " }
    @{ Cell = "D30"; Value = "'170.22" }
    @{ Cell = "E30"; Value = "  -0.44%  " }
    @{ Cell = "D31"; Value = "'20.72" }
    @{ Cell = "E31"; Value = "  +0.51%  " }
    @{ Cell = "E32"; Value = "  -1.33%  " }
    @{ Cell = "D33"; Value = "'5.53" }
    @{ Cell = "E33"; Value = "  +3.52%  " }
    @{ Cell = "D34"; Value = "'0.123" }
    @{ Cell = "E34"; Value = "  -2.83%  " }
    @{ Cell = "D35"; Value = "'0.0738" }
    @{ Cell = "E35"; Value = "  +2.16%  " }
    @{ Cell = "E36"; Value = "  +0.26%  " }
    @{ Cell = "D37"; Value = "'26.66" }
    @{ Cell = "E37"; Value = "  +15.66%  " }
    @{ Cell = "D38"; Value = "'4.02" }
    @{ Cell = "E38"; Value = "  +5.63%  " }
    @{ Cell = "E39"; Value = "  +8.14%  " }
    @{ Cell = "E40"; Value = "  -0.95%  " }
    @{ Cell = "E41"; Value = "  -0.29%  " }
    @{ Cell = "D42"; Value = "'12.43" }
    @{ Cell = "E42"; Value = "  +21.68%  " }
    @{ Cell = "D43"; Value = "'65.55" }
    @{ Cell = "E43"; Value = "  -0.91%  " }
    @{ Cell = "E44"; Value = "  +6.27%  " }
    @{ Cell = "E45"; Value = "  -2.95%  " }
    @{ Cell = "D46"; Value = "'4.75" }
    @{ Cell = "E46"; Value = "  +2.85%  " }
    @{ Cell = "D47"; Value = "'8.66" }
    @{ Cell = "E47"; Value = "  -8.10%  " }
    @{ Cell = "E48"; Value = "  -1.80%  " }
    @{ Cell = "E49"; Value = "  +0.17%  " }
    @{ Cell = "E50"; Value = "  +4.32%  " }
    @{ Cell = "E51"; Value = "  -0.32%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

